$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect to make the edits, then
# restore protection afterward.
$ws.Unprotect()

# Update the confidential disclaimer date (A12) from 2021-04-08 to 2021-04-09
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."
$ws.Rows(12).AutoFit()

# Update the Weight (D) and Percent Change (E) values for rows 2-9
$ws.Range("D2").Value = 0.1776204984956161
$ws.Range("E2").Value = -0.0009025270758121762

$ws.Range("D3").Value = 0.1774722142346879
$ws.Range("E3").Value = -0.0009823182711198308

$ws.Range("D4").Value = 0.2252057193640208
$ws.Range("E4").Value = -0.002495840266222848

$ws.Range("D5").Value = 0.07993824160916477

$ws.Range("D6").Value = 0.07982502457210473

$ws.Range("D7").Value = 0.1203547199982366
$ws.Range("E7").Value = -0.0009823182711199419

$ws.Range("D8").Value = 0.1395835817261691
$ws.Range("E8").Value = -0.0008326394671107629

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = -0.001131168449918452

# Restore sheet protection
$ws.Protect()
